$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = '@'
$ws.Cells.Item(2, 4).Value = '72.432.54'
$ws.Cells.Item(2, 5).Value = '  +4.76%  '

$ws.Cells.Item(3, 4).NumberFormat = '@'
$ws.Cells.Item(3, 4).Value = '4.050.74'
$ws.Cells.Item(3, 5).Value = '  +4.19%  '

$ws.Cells.Item(4, 5).Value = '  -0.09%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '521.66'
$ws.Cells.Item(5, 5).Value = '  -0.84%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '147.99'
$ws.Cells.Item(6, 5).Value = '  +2.42%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.713'
$ws.Cells.Item(7, 5).Value = '  +16.46%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '4.040.49'
$ws.Cells.Item(8, 5).Value = '  +4.18%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.999'
$ws.Cells.Item(9, 5).Value = '  +0.10%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.772'
$ws.Cells.Item(10, 5).Value = '  +7.81%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.180'
$ws.Cells.Item(11, 5).Value = '  +5.53%  '

$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.0000335'
$ws.Cells.Item(12, 5).Value = '  +1.31%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '48.47'
$ws.Cells.Item(13, 5).Value = '  +15.72%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '11.18'
$ws.Cells.Item(14, 5).Value = '  +9.70%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '4.681.99'
$ws.Cells.Item(15, 5).Value = '  +3.63%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '4.058.16'
$ws.Cells.Item(16, 5).Value = '  +4.37%  '

$ws.Cells.Item(17, 4).NumberFormat = '@'
$ws.Cells.Item(17, 4).Value = '21.33'
$ws.Cells.Item(17, 5).Value = '  +8.23%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '14.30'
$ws.Cells.Item(18, 5).Value = '  +2.47%  '

$ws.Cells.Item(19, 5).Value = '  +0.98%  '

$ws.Cells.Item(20, 5).Value = '  -0.32%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '72.292.35'
$ws.Cells.Item(21, 5).Value = '  +4.47%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '446.04'
$ws.Cells.Item(22, 5).Value = '  +5.67%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '105.03'
$ws.Cells.Item(23, 5).Value = '  +19.97%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '3.61'
$ws.Cells.Item(24, 5).Value = '  +6.52%  '

$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '15.26'
$ws.Cells.Item(25, 5).Value = '  +8.19%  '

$ws.Cells.Item(26, 5).Value = '  +1.53%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '11.51'
$ws.Cells.Item(27, 5).Value = '  +0.96%  '

$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '11.12'
$ws.Cells.Item(28, 5).Value = '  +5.50%  '

$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '38.02'
$ws.Cells.Item(29, 5).Value = '  +4.96%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '5.82'
$ws.Cells.Item(30, 5).Value = '  +2.41%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '3.32'
$ws.Cells.Item(31, 5).Value = '  +16.55%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '13.81'
$ws.Cells.Item(32, 5).Value = '  +4.92%  '

$ws.Cells.Item(33, 5).Value = '  +3.95%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '682.16'
$ws.Cells.Item(34, 5).Value = '  -1.38%  '

$ws.Cells.Item(35, 5).Value = '  +14.53%  '

$ws.Cells.Item(36, 5).Value = '  -0.57%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '42.56'
$ws.Cells.Item(37, 5).Value = '  +6.80%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '0.0₃0874'
$ws.Cells.Item(38, 5).Value = '  +2.70%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.431'
$ws.Cells.Item(39, 5).Value = '  +0.50%  '

$ws.Cells.Item(40, 2).Value = 'Kaspa'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '0.153'
$ws.Cells.Item(40, 5).Value = '  +3.63%  '

$ws.Cells.Item(41, 2).Value = 'ThetaToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '3.49'
$ws.Cells.Item(41, 5).Value = '  +6.18%  '

$ws.Cells.Item(42, 5).Value = '  +0.07%  '

$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '0.0503'
$ws.Cells.Item(43, 5).Value = '  +4.79%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.997'
$ws.Cells.Item(44, 5).Value = '  -0.38%  '

$ws.Cells.Item(45, 5).Value = '  -1.08%  '

$ws.Cells.Item(46, 5).Value = '  +11.74%  '

$ws.Cells.Item(47, 2).Value = 'Fetch.AI'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '2.72'
$ws.Cells.Item(47, 5).Value = '  -1.45%  '

$ws.Cells.Item(48, 2).Value = 'ApeXProtocol'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '3.49'
$ws.Cells.Item(48, 5).Value = '  +2.39%  '

$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '9.57'
$ws.Cells.Item(49, 5).Value = '  +12.25%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '3.09'
$ws.Cells.Item(50, 5).Value = '  +3.58%  '

$ws.Cells.Item(51, 2).Value = 'FLOKI'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Cells.Item(51, 4).NumberFormat = '@'
$ws.Cells.Item(51, 4).Value = '0.000275'
$ws.Cells.Item(51, 5).Value = '  +3.35%  '
